$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.00560353008272598
$ws.Range("E2").Value = 0.00560353008272598

# Row 3
$ws.Range("D3").Value = 0.04049788347242873
$ws.Range("E3").Value = 0.04049788347242873

# Row 4
$ws.Range("D4").Value = 0.003162751500676715
$ws.Range("E4").Value = 0.003162751500676715

# Row 5
$ws.Range("D5").Value = 0.005287047241691047
$ws.Range("E5").Value = 0.005287047241691047

# Row 6
$ws.Range("D6").Value = 0.07402200682465487
$ws.Range("E6").Value = 0.07402200682465487

# Row 7
$ws.Range("D7").Value = 0.9879630327676151
$ws.Range("E7").Value = 0.0120369672323849

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.05128823871843159
$ws.Range("E8").Value = 0.9487117612815684

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.05288793453929659
$ws.Range("E9").Value = 0.9471120654607034

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.09847752772220425
$ws.Range("E10").Value = 0.9015224722777957

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.1177659684766499
$ws.Range("E11").Value = 0.8822340315233501
$ws.Range("F11").Value = 1.051130056381226
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.003089965693185638
$ws.Range("E12").Value = 0.003089965693185638

# Row 13
$ws.Range("D13").Value = [double]"9.749542322470059E-05"
$ws.Range("E13").Value = [double]"9.749542322470059E-05"

# Row 14
$ws.Range("D14").Value = 0.001431768515885918
$ws.Range("E14").Value = 0.001431768515885918

# Row 15
$ws.Range("D15").Value = 0.0001870784816207065
$ws.Range("E15").Value = 0.0001870784816207065

# Row 16
$ws.Range("D16").Value = 0.05683324498390609
$ws.Range("E16").Value = 0.05683324498390609

# Row 17
$ws.Range("D17").Value = 0.998441086860117
$ws.Range("E17").Value = 0.001558913139882989

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.003774482659735072
$ws.Range("E18").Value = 0.9962255173402649

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"6.442129377297509E-06"
$ws.Range("E19").Value = 0.9999935578706227

# Row 20
$ws.Range("D20").Value = 0.6473810642075165
$ws.Range("E20").Value = 0.3526189357924835

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.01766068135839209
$ws.Range("E21").Value = 0.982339318641608
$ws.Range("F21").Value = 2.206826448440552
$ws.Range("G21").Value = 0.7
